# This workbook is a weekly price log. A new week of data is being added
# for "Cebollín baby" (Agrícola del Norte S.A. de Arica): a new record is
# inserted at row 123, pushing all the existing records (previously rows
# 123-141) down by one row (now rows 124-142).
#
# The new row 123 keeps the same fixed attributes as the record that used
# to occupy that row (now at row 124) - market/region/category/etc. - but
# gets its own date (column D) and price figures (columns K, L, M, P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 123; everything below (123-141) shifts to 124-142.
$ws.Rows.Item(123).Insert()

# Seed the new row 123 with the same values as the row right below it
# (which now holds what used to be row 123), then overwrite the cells
# that actually carry new data for this week's entry.
$lastCol = 18
for ($c = 1; $c -le $lastCol; $c++) {
    $ws.Cells.Item(123, $c).Value = $ws.Cells.Item(124, $c).Value()
}

$ws.Cells.Item(123, 4).Value  = 45180   # D123 - Fecha
$ws.Cells.Item(123, 11).Value = 900     # K123 - Precio mínimo
$ws.Cells.Item(123, 12).Value = 1000    # L123 - Precio máximo
$ws.Cells.Item(123, 13).Value = 950     # M123 - Precio promedio ponderado
$ws.Cells.Item(123, 16).Value = 475     # P123 - Precio $/Kg
